$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "1.00", "6.90")
# must be forced to Text format first, otherwise Excel auto-converts them to
# numbers and loses the original formatting (trailing zeros, etc).
$textCells = @("D4", "D5", "D6", "D11", "D18", "D19", "D20", "D21", "D22", "D24", "D29", "D31", "D32", "D34", "D35", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto data refresh.
$ws.Range("D2").Value = "66.774.09"
$ws.Range("E2").Value = "  -4.74%  "
$ws.Range("D3").Value = "3.469.72"
$ws.Range("E3").Value = "  -6.13%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "604.25"
$ws.Range("E5").Value = "  -7.43%  "
$ws.Range("D6").Value = "147.72"
$ws.Range("E6").Value = "  -9.25%  "
$ws.Range("D7").Value = "3.468.16"
$ws.Range("E7").Value = "  -6.08%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -4.31%  "
$ws.Range("E10").Value = "  -6.54%  "
$ws.Range("D11").Value = "6.90"
$ws.Range("E11").Value = "  -4.28%  "
$ws.Range("E12").Value = "  -5.78%  "
$ws.Range("E13").Value = "  -8.31%  "
$ws.Range("D14").Value = "4.053.81"
$ws.Range("E14").Value = "  -6.04%  "
$ws.Range("E15").Value = "  -5.60%  "
$ws.Range("D16").Value = "3.476.18"
$ws.Range("E16").Value = "  -5.66%  "
$ws.Range("D17").Value = "66.674.78"
$ws.Range("E17").Value = "  -4.69%  "
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").Value = "14.91"
$ws.Range("E20").Value = "  -7.24%  "
$ws.Range("D21").Value = "440.86"
$ws.Range("E21").Value = "  -6.75%  "
$ws.Range("D22").Value = "8.95"
$ws.Range("E22").Value = "  -14.05%  "
$ws.Range("E23").Value = "  -5.20%  "
$ws.Range("D24").Value = "76.99"
$ws.Range("E24").Value = "  -3.72%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "3.606.28"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("E28").Value = "  -10.02%  "
$ws.Range("D29").Value = "8.28"
$ws.Range("E29").Value = "  -8.18%  "
$ws.Range("E30").Value = "  -5.86%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.56"
$ws.Range("E31").Value = "  -9.62%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("D34").Value = "25.50"
$ws.Range("E34").Value = "  -5.04%  "
$ws.Range("D35").Value = "6.10"
$ws.Range("E35").Value = "  -7.32%  "
$ws.Range("E36").Value = "  -8.55%  "
$ws.Range("D37").Value = "3.451.92"
$ws.Range("E37").Value = "  -6.37%  "
$ws.Range("D38").Value = "7.89"
$ws.Range("E38").Value = "  -6.83%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.17"
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "172.18"
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("D43").Value = "0.0865"
$ws.Range("E44").Value = "  -8.69%  "
$ws.Range("D45").Value = "0.878"
$ws.Range("E45").Value = "  -5.98%  "
$ws.Range("D46").Value = "45.62"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").Value = "25.91"
$ws.Range("E48").Value = "  -12.27%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "2.48"
$ws.Range("E49").Value = "  -13.50%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "7.53"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -6.11%  "
